$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Client")

# Insert a new column before column K ("Mobile"), shifting the "Mobile".."Status"
# header block (and everything below it) one column to the right, so a new
# column for "White Label" appears right after "Limit".
$ws.Columns("K").Insert(-4161)  # -4161 = xlShiftToRight

# Give the new header cell K2 the same look as its left neighbour (J2 = "Limit"),
# then set its text to the new field name.
$ws.Range("J2").Copy()
$ws.Range("K2").PasteSpecial(-4122)  # -4122 = xlPasteFormats
$ws.Range("K2").Value = "White Label"

# Mirror the cursor position left behind in the authored workbook (cosmetic).
$ws.Range("K7").Select()
